$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50, shifting existing rows 50:130 down to 51:131
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new data record
$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "Terminal La Palmera de La Serena"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 44665
$ws.Range("D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 100112001
$ws.Range("G50").Value = "Berenjena"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 400
$ws.Range("K50").Value = 8000
$ws.Range("L50").Value = 9000
$ws.Range("M50").Value = 8500
$ws.Range("N50").Value = "`$/caja 50 unidades"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 170
$ws.Range("Q50").Value = 50
$ws.Range("R50").Value = "Hortaliza"
